# Tyoajanseurannan ja product backlogin paivitys
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jaana")

# Row 14: new time entry.
#  - copy the date-formatted style from A13 onto A14 so the number format matches
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the description text first so the shared-string table gets the same
# ordering as the authored workbook, then the date and hours.
$ws.Range("C14").Value = "CustomerWindow ja Customer-class toimintoja"
$d = Get-Date -Year 2023 -Month 3 -Day 1
$ws.Range("A14").Value = $d.Date
$ws.Range("B14").Value = 5

# The description text wraps onto two lines, so the row grows taller to fit it
$ws.Rows.Item(14).RowHeight = 31.8

# Row 4: add person's name next to "Henkilo:" label
$ws.Range("B4").Value = "Jaana Pusa"

# Update selection to reflect where the user left off editing
$ws.Range("B15").Select()

$wb.Save()
